# Auto-generated edit script applying cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "26.915.27"
$c.Style = "Normal"

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.96%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.819.11"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.74%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "309.43"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4676"
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.10%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3693"
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.07373"
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.07%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.8718"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.34%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "20.46"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.61%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.822.79"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.97%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.368"
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.74%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "92.58"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.55%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.07072"
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.87%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "6.503"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.69%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.60%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "14.77"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.47%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "26.943.17"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.06%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "5.349"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.86%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.36%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.039.62"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.67%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "151.45"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.25%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.78%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "18.41"
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.95%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "5.328"
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.84%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "115.76"
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.04%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.08944"
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.36%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.7702"
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.67%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.36%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "4.503"
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.70%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.05%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.84%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.01967"
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.50%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.05290"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.75%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 2)
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "7.322"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.06%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = "@"
$c.Value = "MXToken"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.950"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.62%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.5344"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.30%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "2.363"
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.72%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.1671"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.03%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "8.453"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.11%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.4960"
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.77%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.54%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 2)
$c.NumberFormat = "@"
$c.Value = "Quant"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "104.18"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.65%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = "PaxDollar"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.672"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  +2.06%  "
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.Style = "Normal"

Write-Host "Applied $([int]97) cell updates"
